# "Natmi following Dr Hou advice"
# Recomputed the Adam23-Itga4 ligand/receptor NATMI edge table: corrected the
# expressing-cell counts/specificity metrics for the existing Sending-cluster x
# Target-cluster combinations (ECs/FAPs) and added the missing "sCs" sending
# cluster rows, completing the full 3x3 cluster combination grid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-10, columns A-T (20 columns). Columns A-D are text labels,
# columns E-T are numeric metrics.
$data = @(
    ,@("ECs", "Adam23", "Itga4", "ECs", 3, 1, 0.1145763333333333, 0.343729, 0.006557053879060051, 0.006557053879060051, 3, 1, 26.04517333333333, 78.13552, 0.9210237118384171, 0.921023711838417, 2.984160461564444, 26.85744415408, 0.00603920210241638, 0.006039202102416379)
    ,@("ECs", "Adam23", "Itga4", "FAPs", 3, 1, 0.1145763333333333, 0.343729, 0.006557053879060051, 0.006557053879060051, 2, 0.6666666666666666, 0.3302223333333333, 0.990667, 0.01167750336256582, 0.01167750336256582, 0.03783566413811111, 0.340520977243, 0.00007657001872124902, 0.00007657001872124902)
    ,@("ECs", "Adam23", "Itga4", "sCs", 3, 1, 0.1145763333333333, 0.343729, 0.006557053879060051, 0.006557053879060051, 3, 1, 1.903109, 5.709327, 0.06729878479901708, 0.06729878479901708, 0.2180512511536666, 1.962461260383, 0.0004412817579224225, 0.0004412817579224225)
    ,@("FAPs", "Adam23", "Itga4", "ECs", 3, 1, 9.390663666666667, 28.171991, 0.5374154140831726, 0.5374154140831726, 3, 1, 26.04517333333333, 78.13552, 0.9210237118384171, 0.921023711838417, 244.5814629133689, 2201.23316622032, 0.4949723394780636, 0.4949723394780635)
    ,@("FAPs", "Adam23", "Itga4", "FAPs", 3, 1, 9.390663666666667, 28.171991, 0.5374154140831726, 0.5374154140831726, 2, 0.6666666666666666, 0.3302223333333333, 0.990667, 0.01167750336256582, 0.01167750336256582, 3.101006867555222, 27.909061807997, 0.006275670305050953, 0.006275670305050953)
    ,@("FAPs", "Adam23", "Itga4", "sCs", 3, 1, 9.390663666666667, 28.171991, 0.5374154140831726, 0.5374154140831726, 3, 1, 1.903109, 5.709327, 0.06729878479901708, 0.06729878479901708, 17.87145654000633, 160.843108860057, 0.03616740430005808, 0.03616740430005808)
    ,@("sCs", "Adam23", "Itga4", "ECs", 3, 1, 7.968512, 23.905536, 0.4560275320377672, 0.4560275320377672, 3, 1, 26.04517333333333, 78.13552, 0.9210237118384171, 0.921023711838417, 207.5412762487467, 1867.87148623872, 0.4200121702579371, 0.420012170257937)
    ,@("sCs", "Adam23", "Itga4", "FAPs", 3, 1, 7.968512, 23.905536, 0.4560275320377672, 0.4560275320377672, 2, 0.6666666666666666, 0.3302223333333333, 0.990667, 0.01167750336256582, 0.01167750336256582, 2.631380625834666, 23.682425632512, 0.00532526303879362, 0.00532526303879362)
    ,@("sCs", "Adam23", "Itga4", "sCs", 3, 1, 7.968512, 23.905536, 0.4560275320377672, 0.4560275320377672, 3, 1, 1.903109, 5.709327, 0.06729878479901708, 0.06729878479901708, 15.164946903808, 136.484522134272, 0.03069009874103656, 0.03069009874103656)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}

Write-Output "applied"
